$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.433.71'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.894.60'
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.21'
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.30'
$ws.Range("E6").Value = '  -5.95%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.548'
$ws.Range("E8").Value = '  -3.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.898.05'
$ws.Range("E9").Value = '  -2.71%  '
$ws.Range("E10").Value = '  -5.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.98'
$ws.Range("E11").Value = '  -2.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.358'
$ws.Range("E12").Value = '  -2.55%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.398.78'
$ws.Range("E13").Value = '  -2.68%  '
$ws.Range("E14").Value = '  +2.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.466.86'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.66'
$ws.Range("E16").Value = '  -3.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.901.63'
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("E18").Value = '  -4.02%  '
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.61'
$ws.Range("E20").Value = '  -3.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.44'
$ws.Range("E21").Value = '  -7.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.56'
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.70'
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.47'
$ws.Range("E25").Value = '  -1.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.453'
$ws.Range("E26").Value = '  -3.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.179'
$ws.Range("E27").Value = '  -5.54%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.81'
$ws.Range("E29").Value = '  -4.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0835'
$ws.Range("E30").Value = '  -10.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.55'
$ws.Range("E33").Value = '  -4.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.92'
$ws.Range("E34").Value = '  -6.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.31'
$ws.Range("E35").Value = '  -7.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.56'
$ws.Range("E36").Value = '  -5.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.994'
$ws.Range("E37").Value = '  -7.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.20'
$ws.Range("E38").Value = '  -4.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.64'
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("E40").Value = '  -4.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.71'
$ws.Range("E41").Value = '  -5.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.283.42'
$ws.Range("E42").Value = '  -5.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.645'
$ws.Range("E43").Value = '  -3.38%  '
$ws.Range("E44").Value = '  -1.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.26'
$ws.Range("E45").Value = '  -8.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.92'
$ws.Range("E47").Value = '  -3.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0236'
$ws.Range("E48").Value = '  -3.77%  '
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0917'
$ws.Range("E50").Value = '  -3.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '247.43'
$ws.Range("E51").Value = '  -7.10%  '
